# Scheduled market-data refresh.
#
# The leve-profitability sheets (one per crafting class: ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) pull current-market-board average prices into column H
# (currentAveragePrice), split NQ/HQ in I/J (currentAveragePriceNQ/HQ), and
# derive LevePriceNQ/HQ (K/L) and LeveProfitNQ/HQ (M/N) from them. This runner
# writes the latest snapshot values for the rows whose market prices moved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 776.75
$ws.Range("I18").Value = 673.1429000000001
$ws.Range("K18").Value = 673.1429000000001
$ws.Range("M18").Value = -389.1429000000001
$ws.Range("H28").Value = 1464.7
$ws.Range("I28").Value = 1464.7
$ws.Range("K28").Value = 1464.7
$ws.Range("M28").Value = -979.7
$ws.Range("H32").Value = 1441.5294
$ws.Range("I32").Value = 874.5
$ws.Range("J32").Value = 1616
$ws.Range("K32").Value = 874.5
$ws.Range("L32").Value = 1616
$ws.Range("M32").Value = -548.5
$ws.Range("N32").Value = -2268
$ws.Range("H41").Value = 1156.0714
$ws.Range("I41").Value = 496
$ws.Range("J41").Value = 2806.25
$ws.Range("K41").Value = 496
$ws.Range("L41").Value = 2806.25
$ws.Range("M41").Value = -56
$ws.Range("N41").Value = -3686.25
$ws.Range("H57").Value = 109329
$ws.Range("J57").Value = 109329
$ws.Range("L57").Value = 327987
$ws.Range("N57").Value = -328985
$ws.Range("H74").Value = 11747.739
$ws.Range("I74").Value = 11949.95
$ws.Range("K74").Value = 11949.95
$ws.Range("M74").Value = -11013.95
$ws.Range("H77").Value = 11747.739
$ws.Range("I77").Value = 11949.95
$ws.Range("K77").Value = 59749.75
$ws.Range("M77").Value = -55069.75
$ws.Range("H100").Value = 12054.2
$ws.Range("I100").Value = 10943.333
$ws.Range("J100").Value = 12530.286
$ws.Range("K100").Value = 10943.333
$ws.Range("L100").Value = 12530.286
$ws.Range("M100").Value = -10402.333
$ws.Range("N100").Value = -13612.286
$ws.Range("H137").Value = 3244.0908
$ws.Range("I137").Value = 3551.3076
$ws.Range("K137").Value = 10653.9228
$ws.Range("M137").Value = -8103.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2000
$ws.Range("J110").Value = 1000
$ws.Range("L110").Value = 1000
$ws.Range("N110").Value = -5090
$ws.Range("H132").Value = 4050777
$ws.Range("I132").Value = 4527104
$ws.Range("K132").Value = 13581312
$ws.Range("M132").Value = -13578782
$ws.Range("H137").Value = 99247.25
$ws.Range("J137").Value = 99247.25
$ws.Range("L137").Value = 99247.25
$ws.Range("N137").Value = -109447.25
$ws.Range("H139").Value = 104826.6
$ws.Range("J139").Value = 104826.6
$ws.Range("L139").Value = 104826.6
$ws.Range("N139").Value = -115106.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 89442.25
$ws.Range("J59").Value = 128885
$ws.Range("L59").Value = 128885
$ws.Range("N59").Value = -130579
$ws.Range("H64").Value = 1540.1666
$ws.Range("I64").Value = 1506
$ws.Range("J64").Value = 1543.2727
$ws.Range("K64").Value = 1506
$ws.Range("L64").Value = 1543.2727
$ws.Range("M64").Value = -1281
$ws.Range("N64").Value = -1993.2727
$ws.Range("H67").Value = 1540.1666
$ws.Range("I67").Value = 1506
$ws.Range("J67").Value = 1543.2727
$ws.Range("K67").Value = 1506
$ws.Range("L67").Value = 1543.2727
$ws.Range("M67").Value = -726
$ws.Range("N67").Value = -3103.2727
$ws.Range("H86").Value = 3928124.8
$ws.Range("I86").Value = 9528580
$ws.Range("J86").Value = 7806.2
$ws.Range("K86").Value = 9528580
$ws.Range("L86").Value = 7806.2
$ws.Range("M86").Value = -9527457
$ws.Range("N86").Value = -10052.2
$ws.Range("H89").Value = 3928124.8
$ws.Range("I89").Value = 9528580
$ws.Range("J89").Value = 7806.2
$ws.Range("K89").Value = 47642900
$ws.Range("L89").Value = 39031
$ws.Range("M89").Value = -47637284
$ws.Range("N89").Value = -50263
$ws.Range("H105").Value = 4077.5789
$ws.Range("I105").Value = 2765
$ws.Range("K105").Value = 2765
$ws.Range("M105").Value = -1018
$ws.Range("H107").Value = 1795.4736
$ws.Range("I107").Value = 1417.0588
$ws.Range("K107").Value = 1417.0588
$ws.Range("M107").Value = 502.9412

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2779.8572
$ws.Range("I16").Value = 1725
$ws.Range("K16").Value = 1725
$ws.Range("M16").Value = -1438
$ws.Range("H113").Value = 2779.8572
$ws.Range("I113").Value = 1725
$ws.Range("K113").Value = 1725
$ws.Range("M113").Value = 445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7399723.5
$ws.Range("J4").Value = 4523955
$ws.Range("L4").Value = 13571865
$ws.Range("N4").Value = -13572089
$ws.Range("H96").Value = 20004.5
$ws.Range("J96").Value = 20004.5
$ws.Range("L96").Value = 60013.5
$ws.Range("N96").Value = -64131.5
$ws.Range("H117").Value = 372.25
$ws.Range("J117").Value = 464.4
$ws.Range("L117").Value = 1393.2
$ws.Range("N117").Value = -8277.200000000001
$ws.Range("H122").Value = 1311.3334
$ws.Range("I122").Value = 931
$ws.Range("K122").Value = 8379
$ws.Range("M122").Value = -5929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 101998
$ws.Range("J135").Value = 101998
$ws.Range("L135").Value = 101998
$ws.Range("N135").Value = -112138

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1026.5518
$ws.Range("I16").Value = 360.09525
$ws.Range("J16").Value = 2776
$ws.Range("K16").Value = 360.09525
$ws.Range("L16").Value = 2776
$ws.Range("M16").Value = -190.09525
$ws.Range("N16").Value = -3116
$ws.Range("H40").Value = 4295.234
$ws.Range("I40").Value = 4189.7075
$ws.Range("K40").Value = 4189.7075
$ws.Range("M40").Value = -4053.7075
$ws.Range("H61").Value = 1121.5358
$ws.Range("I61").Value = 1150.2307
$ws.Range("J61").Value = 748.5
$ws.Range("K61").Value = 1150.2307
$ws.Range("L61").Value = 748.5
$ws.Range("M61").Value = -948.2307000000001
$ws.Range("N61").Value = -1152.5
$ws.Range("H100").Value = 3975
$ws.Range("I100").Value = 3528.5715
$ws.Range("J100").Value = 4600
$ws.Range("K100").Value = 3528.5715
$ws.Range("L100").Value = 4600
$ws.Range("M100").Value = -2987.5715
$ws.Range("N100").Value = -5682
$ws.Range("H113").Value = 1121.5358
$ws.Range("I113").Value = 1150.2307
$ws.Range("J113").Value = 748.5
$ws.Range("K113").Value = 1150.2307
$ws.Range("L113").Value = 748.5
$ws.Range("M113").Value = 1019.7693
$ws.Range("N113").Value = -5088.5
$ws.Range("H132").Value = 2881.8767
$ws.Range("J132").Value = 2874.9583
$ws.Range("L132").Value = 8624.874899999999
$ws.Range("N132").Value = -13684.8749
$ws.Range("H133").Value = 88281.75
$ws.Range("J133").Value = 88281.75
$ws.Range("L133").Value = 88281.75
$ws.Range("N133").Value = -93341.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 898332.7
$ws.Range("I2").Value = 898332.7
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 898332.7
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -898220.7
$ws.Range("N2").ClearContents()
$ws.Range("H107").Value = 1127.8422
$ws.Range("I107").Value = 912
$ws.Range("J107").Value = 1595.5
$ws.Range("K107").Value = 2736
$ws.Range("L107").Value = 4786.5
$ws.Range("M107").Value = -816
$ws.Range("N107").Value = -8626.5
$ws.Range("H132").Value = 2027.5938
$ws.Range("J132").Value = 1644.75
$ws.Range("L132").Value = 4934.25
$ws.Range("N132").Value = -9994.25
$ws.Range("H133").Value = 107362.8
$ws.Range("J133").Value = 113604.664
$ws.Range("L133").Value = 113604.664
